$wb = $excel.ActiveWorkbook

# Sheet 1: JFK
$ws = $wb.Worksheets.Item("JFK")
$ws.Range("C2").Value = 177271
$ws.Range("D2").Value = 1897495.680000036
$ws.Range("E2").Value = 112301

$ws.Range("C3").Value = 60379
$ws.Range("D3").Value = 634690.0400000043
$ws.Range("E3").Value = 37880

$ws.Range("C4").Value = 149863
$ws.Range("D4").Value = 1628945.43000004
$ws.Range("E4").Value = 96089

$ws.Range("C5").Value = 60561
$ws.Range("D5").Value = 632690.6500000041
$ws.Range("E5").Value = 37912

$ws.Range("C6").Value = 60448
$ws.Range("D6").Value = 660756.8200000003
$ws.Range("E6").Value = 39258

$ws.Range("C7").Value = 25681
$ws.Range("D7").Value = 275461.1300000007
$ws.Range("E7").Value = 16500

# Sheet 2: Regular
$ws = $wb.Worksheets.Item("Regular")
$ws.Range("C2").Value = 6897913
$ws.Range("D2").Value = 11369332.68999879
$ws.Range("E2").Value = 4618013

$ws.Range("C3").Value = 2308458
$ws.Range("D3").Value = 3642218.109999947
$ws.Range("E3").Value = 1451048

$ws.Range("C4").Value = 6335233
$ws.Range("D4").Value = 10423945.63999971
$ws.Range("E4").Value = 4291341

$ws.Range("C5").Value = 2706845
$ws.Range("D5").Value = 4264982.189999899
$ws.Range("E5").Value = 1713625

$ws.Range("C6").Value = 3112053
$ws.Range("D6").Value = 5360842.019999959
$ws.Range("E6").Value = 2153295

$ws.Range("C7").Value = 1085008
$ws.Range("D7").Value = 1856013.899999982
$ws.Range("E7").Value = 703879

# Sheet 3: Others
$ws = $wb.Worksheets.Item("Others")
$ws.Range("C2").Value = 38527
$ws.Range("D2").Value = 1004739.729999992
$ws.Range("E2").Value = 77782

$ws.Range("C3").Value = 13398
$ws.Range("D3").Value = 196196.8399999999
$ws.Range("E3").Value = 19272

$ws.Range("C4").Value = 33635
$ws.Range("D4").Value = 844561.2199999882
$ws.Range("E4").Value = 60647

$ws.Range("C5").Value = 14653
$ws.Range("D5").Value = 206438.879999999
$ws.Range("E5").Value = 20087

$ws.Range("C6").Value = 14896
$ws.Range("D6").Value = 1046400.299999988
$ws.Range("E6").Value = 40429

$ws.Range("C7").Value = 6690
$ws.Range("D7").Value = 146991.6700000001
$ws.Range("E7").Value = 11522
